$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.2453501895176297
$ws.Range("C2").Value = 1.885822625406119
$ws.Range("D2").Value = 16.01396840704841
$ws.Range("E2").Value = 4.0017456699606
$ws.Range("F2").Value = 4.083986181670692
$ws.Range("G2").Value = 23
$ws.Range("B3").Value = -0.008306120572730035
$ws.Range("C3").Value = 1.678411150852771
$ws.Range("D3").Value = 10.3892641846745
$ws.Range("E3").Value = 3.22323815202577
$ws.Range("F3").Value = 3.299078473400757
$ws.Range("G3").Value = 22
$ws.Range("B4").Value = -0.6206486848402303
$ws.Range("C4").Value = 1.081319350267402
$ws.Range("D4").Value = 4.95937692841922
$ws.Range("E4").Value = 2.226965857039398
$ws.Range("F4").Value = 2.191547568579468
$ws.Range("G4").Value = 21
$ws.Range("B5").Value = -0.007917461514596722
$ws.Range("C5").Value = 0.7359397285945519
$ws.Range("D5").Value = 2.052934180666193
$ws.Range("E5").Value = 1.432806400274019
$ws.Range("F5").Value = 1.470005905634166
$ws.Range("G5").Value = 20
$ws.Range("B6").Value = -0.01396789604063663
$ws.Range("C6").Value = 0.6794082136444377
$ws.Range("D6").Value = 0.985853761070935
$ws.Range("E6").Value = 0.992901687515403
$ws.Range("F6").Value = 1.020008565325462
$ws.Range("G6").Value = 19
$ws.Range("B7").Value = 0.006227155691902452
$ws.Range("C7").Value = 0.5859396340579296
$ws.Range("D7").Value = 0.8714568099673092
$ws.Range("E7").Value = 0.9335185107791432
$ws.Range("F7").Value = 0.9605612509173802
$ws.Range("G7").Value = 18
$ws.Range("B8").Value = 0.1323879842408515
$ws.Range("C8").Value = 0.5124844762799831
$ws.Range("D8").Value = 0.586790449754057
$ws.Range("E8").Value = 0.7660224864545799
$ws.Range("F8").Value = 0.7777164414773035
$ws.Range("G8").Value = 17
$ws.Range("B9").Value = 0.1676475329133778
$ws.Range("C9").Value = 0.4867046202604076
$ws.Range("D9").Value = 0.3917272528104276
$ws.Range("E9").Value = 0.6258811810642877
$ws.Range("F9").Value = 0.6227864759983028
$ws.Range("G9").Value = 16
$ws.Range("B10").Value = 0.199895990469669
$ws.Range("C10").Value = 0.4747094579524145
$ws.Range("D10").Value = 0.3665317847062111
$ws.Range("E10").Value = 0.6054186854617316
$ws.Range("F10").Value = 0.5915234970278875
$ws.Range("G10").Value = 15
$ws.Range("B11").Value = 0.263452728002003
$ws.Range("C11").Value = 0.3988499978544717
$ws.Range("D11").Value = 0.2138270447762498
$ws.Range("E11").Value = 0.4624143648030949
$ws.Range("F11").Value = 0.3943715417630878
$ws.Range("G11").Value = 14
